# Insert a new row at position 89, shifting existing rows 89:160 down to 90:161,
# and populate the new row 89 with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 89 (pushes old rows 89.. down by one).
$ws.Rows("89:89").Insert()

# Fill in the new row 89 with the new data record.
$ws.Cells.Item(89, 1).Value = 9
$ws.Cells.Item(89, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(89, 3).Value = "Metropolitana"
$ws.Cells.Item(89, 4).Value = 44512
$ws.Cells.Item(89, 5).Value = 13
$ws.Cells.Item(89, 6).Value = 100112030
$ws.Cells.Item(89, 7).Value = "Poroto granado"
$ws.Cells.Item(89, 8).Value = "Sin especificar"
$ws.Cells.Item(89, 9).Value = "Primera"
$ws.Cells.Item(89, 10).Value = 34
$ws.Cells.Item(89, 11).Value = 32000
$ws.Cells.Item(89, 12).Value = 34000
$ws.Cells.Item(89, 13).Value = 33000
$ws.Cells.Item(89, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(89, 15).Value = "Perú"
$ws.Cells.Item(89, 16).Value = 1320
$ws.Cells.Item(89, 17).Value = 25
$ws.Cells.Item(89, 18).Value = "Hortaliza"
